# Lattice multiplication exercises: regenerate the 15 practice problems
# (each table cell holds a "AA x BB" problem plus its lattice grid digits).
# We rewrite the problem/digit text per cell while keeping every other
# part of the OOXML (run properties, breaks, table/cell formatting, the
# constant "  ----" divider line, xml:space handling, etc.) untouched.

$d = $word.ActiveDocument

# Map: old top line -> new top line, old 2nd line -> new 2nd line,
#      old 4th line -> new 4th line, old 5th line -> new 5th line.
# (3rd line is always the unchanged "  ----" separator.)
$cellEdits = @(
    ,@("90 x 16", "73 x 34", "  1    6", "  3    4", "9|    |", "7|    |", "0|    |", "3|    |")
    ,@("50 x 11", "56 x 92", "  1    1", "  9    2", "5|    |", "5|    |", "0|    |", "6|    |")
    ,@("80 x 28", "73 x 92", "  2    8", "  9    2", "8|    |", "7|    |", "0|    |", "3|    |")
    ,@("87 x 97", "59 x 21", "  9    7", "  2    1", "8|    |", "5|    |", "7|    |", "9|    |")
    ,@("80 x 63", "72 x 11", "  6    3", "  1    1", "8|    |", "7|    |", "0|    |", "2|    |")
    ,@("84 x 46", "18 x 32", "  4    6", "  3    2", "8|    |", "1|    |", "4|    |", "8|    |")
    ,@("27 x 93", "58 x 36", "  9    3", "  3    6", "2|    |", "5|    |", "7|    |", "8|    |")
    ,@("45 x 96", "50 x 97", "  9    6", "  9    7", "4|    |", "5|    |", "5|    |", "0|    |")
    ,@("77 x 23", "75 x 38", "  2    3", "  3    8", "7|    |", "7|    |", "7|    |", "5|    |")
    ,@("58 x 95", "35 x 51", "  9    5", "  5    1", "5|    |", "3|    |", "8|    |", "5|    |")
    ,@("37 x 84", "66 x 47", "  8    4", "  4    7", "3|    |", "6|    |", "7|    |", "6|    |")
    ,@("42 x 30", "58 x 47", "  3    0", "  4    7", "4|    |", "5|    |", "2|    |", "8|    |")
    ,@("74 x 11", "58 x 16", "  1    1", "  1    6", "7|    |", "5|    |", "4|    |", "8|    |")
    ,@("85 x 68", "58 x 71", "  6    8", "  7    1", "8|    |", "5|    |", "5|    |", "8|    |")
    ,@("60 x 23", "30 x 46", "  2    3", "  4    6", "6|    |", "3|    |", "0|    |", "0|    |")
)

# Pull the canonical OOXML for the whole document so we can do precise,
# per-cell text-node surgery (keeps xml:space="preserve" exactly where it
# already belongs instead of Range.Text's whole-run rewrite behavior).
$xml = $d.Content.WordOpenXML

$searchFrom = 0
foreach ($edit in $cellEdits) {
    $oldTop = $edit[0]; $newTop = $edit[1]
    $oldL2  = $edit[2]; $newL2  = $edit[3]
    $oldL4  = $edit[4]; $newL4  = $edit[5]
    $oldL5  = $edit[6]; $newL5  = $edit[7]

    $anchor = "<w:t>" + $oldTop + "</w:t>"
    $idx = $xml.IndexOf($anchor, $searchFrom)
    if ($idx -lt 0) {
        throw "Could not find cell anchor: $oldTop"
    }
    $endIdx = $xml.IndexOf("</w:tc>", $idx)
    if ($endIdx -lt 0) {
        throw "Could not find closing </w:tc> after: $oldTop"
    }

    $chunk = $xml.Substring($idx, $endIdx - $idx)

    $nodeMatches = [regex]::Matches($chunk, '<w:t( xml:space="preserve")?>([^<]*)</w:t>')
    if ($nodeMatches.Count -ne 5) {
        throw "Expected 5 text nodes in cell for '$oldTop', found $($nodeMatches.Count)"
    }

    $newValues = @($newTop, $newL2, "  ----", $newL4, $newL5)
    $oldValues = @($oldTop, $oldL2, "  ----", $oldL4, $oldL5)

    # Splice back-to-front so earlier replacements don't shift later offsets.
    $newChunk = $chunk
    for ($i = 4; $i -ge 0; $i--) {
        $m = $nodeMatches[$i]
        if ($m.Groups[2].Value -ne $oldValues[$i]) {
            throw "Unexpected text node content '$($m.Groups[2].Value)' (expected '$($oldValues[$i])') in cell for '$oldTop'"
        }
        $preserveAttr = $m.Groups[1].Value
        $replacement = "<w:t" + $preserveAttr + ">" + $newValues[$i] + "</w:t>"
        $newChunk = $newChunk.Substring(0, $m.Index) + $replacement + $newChunk.Substring($m.Index + $m.Length)
    }

    $xml = $xml.Substring(0, $idx) + $newChunk + $xml.Substring($endIdx)

    # Next cell search should start after this cell's original anchor
    # position (text lengths may have changed, so resume from idx+1).
    $searchFrom = $idx + 1
}

$d.Content.InsertXML($xml)
